# Normalize the "Recorded By" (column G) values on the Session Analysis
# Results sheet: wherever "System" appears in the comma-separated list of
# recorders, move it to the front of the list (keeping the remaining
# names in their original relative order).
#
# This mirrors the upstream sync that re-ordered the G-column entries so
# that the automated "System" recorder is always listed first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows on which the "Recorded By" value actually changes order as part of
# this sync (matches the rows touched by the upstream commit).
$targetRows = @(
    2, 3, 5, 6, 7, 8, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26,
    28, 29, 31, 32, 33, 34, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52,
    54, 55, 57, 58, 59, 60, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78,
    80, 81, 82, 83, 84, 85, 86, 90, 92, 93, 94, 96, 99, 101, 106, 107, 108, 109, 110, 111,
    112, 116, 118, 119, 120, 122, 125, 127, 132, 133, 134, 135, 136, 137, 138, 142, 144, 145, 146, 148,
    151, 153
)

foreach ($r in $targetRows) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"

        $hasSystem = $false
        $rest = @()
        foreach ($p in $parts) {
            if ($p.Equals("System")) {
                $hasSystem = $true
            } else {
                $rest += $p
            }
        }

        if ($hasSystem) {
            $newParts = @("System") + $rest
            $newVal = $newParts -join ", "
            if (-not $newVal.Equals($val)) {
                $cell.Value = $newVal
            }
        }
    }
}
